# Fruta / hortaliza, semanal
# Insert a new weekly price record for "Haba" (Macroferia Regional de Talca)
# above the existing row 56, pushing the subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 56 (existing rows 56-58 shift down to 57-59)
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with this week's data
$ws.Cells.Item(56, 1).Value = 5
$ws.Cells.Item(56, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(56, 3).Value = "Maule"
$ws.Cells.Item(56, 4).Value = 44516
$ws.Cells.Item(56, 5).Value = 7
$ws.Cells.Item(56, 6).Value = 100112026
$ws.Cells.Item(56, 7).Value = "Haba"
$ws.Cells.Item(56, 8).Value = "Sin especificar"
$ws.Cells.Item(56, 9).Value = "Primera"
$ws.Cells.Item(56, 10).Value = 600
$ws.Cells.Item(56, 11).Value = 5000
$ws.Cells.Item(56, 12).Value = 5000
$ws.Cells.Item(56, 13).Value = 5000
$ws.Cells.Item(56, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(56, 15).Value = "Región del Maule"
$ws.Cells.Item(56, 16).Value = 200
$ws.Cells.Item(56, 17).Value = 25
$ws.Cells.Item(56, 18).Value = "Hortaliza"
